$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2171.6875
$ws.Range("J40").Value = 2204.8
$ws.Range("L40").Value = 2204.8
$ws.Range("N40").Value = -2554.8
$ws.Range("H113").Value = 3811.875
$ws.Range("I113").Value = 2624.75
$ws.Range("J113").Value = 4999
$ws.Range("K113").Value = 2624.75
$ws.Range("L113").Value = 4999
$ws.Range("M113").Value = 629.25
$ws.Range("N113").Value = -11507
$ws.Range("H116").Value = 4863.5454
$ws.Range("I116").Value = 3083.3333
$ws.Range("K116").Value = 3083.3333
$ws.Range("M116").Value = 358.6667000000002
$ws.Range("H132").Value = 870.63336
$ws.Range("I132").Value = 900.2857
$ws.Range("K132").Value = 2700.8571
$ws.Range("M132").Value = -170.8571000000002
$ws.Range("H137").Value = 1562.4445
$ws.Range("I137").Value = 1343.6666
$ws.Range("K137").Value = 4030.9998
$ws.Range("M137").Value = -1480.9998
$ws.Range("H138").Value = 4756.75
$ws.Range("I138").Value = 3678.4666
$ws.Range("K138").Value = 11035.3998
$ws.Range("M138").Value = -5895.399800000001
$ws.Range("H141").Value = 1782.4348
$ws.Range("I141").Value = 1711.9546
$ws.Range("K141").Value = 5135.8638
$ws.Range("M141").Value = 44.13619999999992
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5133.143
$ws.Range("I32").Value = 3042.0527
$ws.Range("K32").Value = 3042.0527
$ws.Range("M32").Value = -2755.0527
$ws.Range("H61").Value = 1874.8334
$ws.Range("I61").Value = 1819.8
$ws.Range("J61").Value = 2150
$ws.Range("K61").Value = 1819.8
$ws.Range("L61").Value = 2150
$ws.Range("M61").Value = -1607.8
$ws.Range("N61").Value = -2574
$ws.Range("H136").Value = 1874.8334
$ws.Range("I136").Value = 1819.8
$ws.Range("J136").Value = 2150
$ws.Range("K136").Value = 5459.4
$ws.Range("L136").Value = 6450
$ws.Range("M136").Value = -2909.4
$ws.Range("N136").Value = -11550
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 3118.0715
$ws.Range("I105").Value = 3444.125
$ws.Range("K105").Value = 3444.125
$ws.Range("M105").Value = -1697.125
$ws.Range("H134").Value = 1321.1428
$ws.Range("I134").Value = 1333.037
$ws.Range("K134").Value = 3999.111
$ws.Range("M134").Value = -1464.111
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4445.0713
$ws.Range("I31").Value = 4260.0835
$ws.Range("K31").Value = 4260.0835
$ws.Range("M31").Value = -3965.0835
$ws.Range("H34").Value = 4445.0713
$ws.Range("I34").Value = 4260.0835
$ws.Range("K34").Value = 4260.0835
$ws.Range("M34").Value = -4058.0835
$ws.Range("H58").Value = 2107.8696
$ws.Range("I58").Value = 1287.6316
$ws.Range("J58").Value = 6004
$ws.Range("K58").Value = 1287.6316
$ws.Range("L58").Value = 6004
$ws.Range("M58").Value = -1084.6316
$ws.Range("N58").Value = -6410
$ws.Range("H132").Value = 1472.8
$ws.Range("I132").Value = 1470.4783
$ws.Range("J132").Value = 1499.5
$ws.Range("K132").Value = 4411.4349
$ws.Range("L132").Value = 4498.5
$ws.Range("M132").Value = -1881.4349
$ws.Range("N132").Value = -9558.5
$ws.Range("H134").Value = 2614.7144
$ws.Range("I134").Value = 2499.4285
$ws.Range("K134").Value = 7498.2855
$ws.Range("M134").Value = -4963.2855
$ws.Range("H136").Value = 2107.8696
$ws.Range("I136").Value = 1287.6316
$ws.Range("J136").Value = 6004
$ws.Range("K136").Value = 3862.8948
$ws.Range("L136").Value = 18012
$ws.Range("M136").Value = -1312.8948
$ws.Range("N136").Value = -23112
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 265.58066
$ws.Range("I12").Value = 237.63637
$ws.Range("J12").Value = 280.95
$ws.Range("K12").Value = 712.9091100000001
$ws.Range("L12").Value = 842.8499999999999
$ws.Range("M12").Value = -539.9091100000001
$ws.Range("N12").Value = -1188.85
$ws.Range("H33").Value = 192.8
$ws.Range("I33").Value = 203.5
$ws.Range("K33").Value = 1221
$ws.Range("M33").Value = -938
$ws.Range("H37").Value = 0
$ws.Range("J37").Value = 0
$ws.Range("L37").Value = 0
$ws.Range("N37").ClearContents() | Out-Null
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H18").Value = 5500002
$ws.Range("I18").Value = 5500002
$ws.Range("K18").Value = 5500002
$ws.Range("M18").Value = -5499709
$ws.Range("H80").Value = 9233.471
$ws.Range("J80").Value = 13772.777
$ws.Range("L80").Value = 13772.777
$ws.Range("N80").Value = -15768.777
$ws.Range("H83").Value = 9233.471
$ws.Range("J83").Value = 13772.777
$ws.Range("L83").Value = 68863.88499999999
$ws.Range("N83").Value = -78847.88499999999
$ws.Range("H96").Value = 20000
$ws.Range("J96").Value = 20000
$ws.Range("L96").Value = 20000
$ws.Range("M96").Value = -25492
$ws.Range("H132").Value = 2757.92
$ws.Range("I132").Value = 2523.9473
$ws.Range("K132").Value = 7571.841899999999
$ws.Range("M132").Value = -5041.841899999999
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 9220
$ws.Range("I16").Value = 9220
$ws.Range("K16").Value = 9220
$ws.Range("M16").Value = -9050
$ws.Range("H22").Value = 2443.8235
$ws.Range("I22").Value = 1741.0344
$ws.Range("K22").Value = 1741.0344
$ws.Range("M22").Value = -1446.0344
$ws.Range("H27").Value = 2443.8235
$ws.Range("I27").Value = 1741.0344
$ws.Range("K27").Value = 1741.0344
$ws.Range("M27").Value = -1634.0344
$ws.Range("H55").Value = 693
$ws.Range("I55").Value = 508.25
$ws.Range("J55").Value = 1247.25
$ws.Range("K55").Value = 508.25
$ws.Range("L55").Value = 1247.25
$ws.Range("M55").Value = -335.25
$ws.Range("N55").Value = -1593.25
$ws.Range("H68").Value = 4941.2
$ws.Range("J68").Value = 5501.5
$ws.Range("L68").Value = 5501.5
$ws.Range("N68").Value = -6999.5
$ws.Range("H71").Value = 4941.2
$ws.Range("J71").Value = 5501.5
$ws.Range("L71").Value = 27507.5
$ws.Range("N71").Value = -34995.5
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 0
$ws.Range("I96").Value = 0
$ws.Range("K96").Value = 0
$ws.Range("M96").ClearContents() | Out-Null
$ws.Range("H104").Value = 7091.25
$ws.Range("J104").Value = 7091.25
$ws.Range("L104").Value = 7091.25
$ws.Range("N104").Value = -14079.25
$ws.Range("H132").Value = 8096.5713
$ws.Range("I132").Value = 1595
$ws.Range("J132").Value = 12972.75
$ws.Range("K132").Value = 4785
$ws.Range("L132").Value = 38918.25
$ws.Range("M132").Value = -2255
$ws.Range("N132").Value = -43978.25
$ws.Range("H136").Value = 2153.7
$ws.Range("I136").Value = 2059.6667
$ws.Range("K136").Value = 6179.000100000001
$ws.Range("M136").Value = -3629.000100000001
